$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.042.96"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -4.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.556.34"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -4.64%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.86"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.40"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -3.58%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -4.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.76"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -3.23%  "

# Row 13
$ws.Range("E13").Value = "  +6.88%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.945.48"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -5.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.544.64"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -6.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.884"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -4.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.23"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -5.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.035.66"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -4.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -4.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.20"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.82%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.67"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -8.11%  "

# Row 24
$ws.Range("E24").Value = "  -3.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -6.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.09"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -6.38%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("E28").Value = "  -2.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.98"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.88%  "

# Row 30
$ws.Range("E30").Value = "  -3.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.70"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.90%  "

# Row 33
$ws.Range("E33").Value = "  -1.98%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.39"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -10.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.15"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -8.84%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0796"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -4.79%  "

# Row 37
$ws.Range("E37").Value = "  -4.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.64"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +10.39%  "

# Row 39
$ws.Range("E39").Value = "  -3.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.00"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -8.92%  "

# Row 41
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +32.39%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.46%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0311"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -3.81%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.89"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.88%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.112.75"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.50%  "

# Row 46
$ws.Range("E46").Value = "  -0.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.12"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.45"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -8.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.798.07"
$ws.Range("D49").NumberFormat = "General"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.15"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -5.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.05%  "

